# The commit inserts one new daily price record for "Coliflor" (Segunda,
# 2023-06-05 / serial 45082) as the new row 63, pushing the existing rows
# 63-155 down to 64-156 (dimension grows from A1:R155 to A1:R156).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 63; Excel shifts rows 63:155 down to 64:156
# and carries the row-above's number formatting onto the new cells.
$ws.Rows.Item(63).Insert()

# Populate the newly inserted row 63 with the new record.
$ws.Range("A63").Value = 1
$ws.Range("B63").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C63").Value = "Arica y Parinacota"
$ws.Range("D63").Value = 45082
$ws.Range("E63").Value = 15
$ws.Range("F63").Value = 100112008
$ws.Range("G63").Value = "Coliflor"
$ws.Range("H63").Value = "Sin especificar"
$ws.Range("I63").Value = "Segunda"
$ws.Range("J63").Value = 1200
$ws.Range("K63").Value = 800
$ws.Range("L63").Value = 900
$ws.Range("M63").Value = 850
$ws.Range("N63").Value = "`$/unidad"
$ws.Range("O63").Value = "Región de Arica y Parinacota"
$ws.Range("P63").Value = 850
$ws.Range("Q63").Value = 1
$ws.Range("R63").Value = "Hortaliza"
